$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update "Valor Mora" total and "Cant. Periodos" count
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 410820
$ws.Range("F13").Value = 3

# ---------------------------------------------------------------------------
# 2) Make room in the worker table for 2 extra rows (new period 2508).
#    Inserting at 20:21 pushes the footer block (rows 24:25 -> 26:27) down,
#    while rows 16:19 (existing data) stay put for now.
# ---------------------------------------------------------------------------
$ws.Rows("20:21").Insert()

# Copy the "last row" (bottom-border) formatting from the old last data row
# (19) down onto the new last data row (21).
$ws.Range("B19:J19").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)

# Copy the regular data-row formatting (row 16) onto rows 19 and 20, which
# now need to look like ordinary (non-last) data rows.
$ws.Range("B16:J16").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$ws.Range("B20:J20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3) Rewrite the worker/period table contents (re-grouped by worker, with a
#    3rd period "2508" added for each worker).
# ---------------------------------------------------------------------------
$rows = @(
  @{ Row=16; Doc="1047477599"; Name="ESTEFANY PEREZ CEBALLOS";             Periodo="2506"; Mora=80000;  Salario=2000000 },
  @{ Row=17; Doc="1047426490"; Name="LIZZETTE DEL ROSARIO HERRERA PEREZ";  Periodo="2506"; Mora=56940;  Salario=1423500 },
  @{ Row=18; Doc="1047477599"; Name="ESTEFANY PEREZ CEBALLOS";             Periodo="2507"; Mora=80000;  Salario=2000000 },
  @{ Row=19; Doc="1047426490"; Name="LIZZETTE DEL ROSARIO HERRERA PEREZ";  Periodo="2507"; Mora=56940;  Salario=1423500 },
  @{ Row=20; Doc="1047477599"; Name="ESTEFANY PEREZ CEBALLOS";             Periodo="2508"; Mora=80000;  Salario=2000000 },
  @{ Row=21; Doc="1047426490"; Name="LIZZETTE DEL ROSARIO HERRERA PEREZ";  Periodo="2508"; Mora=56940;  Salario=1423500 }
)

foreach ($r in $rows) {
  $n = $r.Row
  $ws.Range("B$n").Value = "CC"
  $ws.Range("C$n").Value = $r.Doc
  $ws.Range("D$n").Value = $r.Name
  $ws.Range("E$n").Value = $r.Periodo
  $ws.Range("F$n").Value = $r.Mora
  $ws.Range("G$n").Value = $r.Salario
}
